$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "57.511.00"
Set-TextValue "E2" "  -0.55%  "
Set-TextValue "D3" "3.098.85"
Set-TextValue "E3" "  +1.23%  "
Set-TextValue "E4" "  -0.03%  "
Set-TextValue "D5" "521.85"
Set-TextValue "E5" "  +1.05%  "
Set-TextValue "D6" "140.97"
Set-TextValue "E6" "  -0.95%  "
Set-TextValue "E7" "  -0.03%  "
Set-TextValue "D8" "3.093.34"
Set-TextValue "E8" "  +1.11%  "
Set-TextValue "E9" "  +0.38%  "
Set-TextValue "D10" "7.19"
Set-TextValue "E10" "  -1.11%  "
Set-TextValue "E11" "  -0.48%  "
Set-TextValue "E12" "  +1.80%  "
Set-TextValue "D13" "3.623.25"
Set-TextValue "E13" "  +1.00%  "
Set-TextValue "E14" "  +1.02%  "
Set-TextValue "D15" "25.79"
Set-TextValue "E15" "  -2.27%  "
Set-TextValue "E16" "  -0.41%  "
Set-TextValue "D17" "57.572.33"
Set-TextValue "E17" "  -0.48%  "
Set-TextValue "D18" "3.093.39"
Set-TextValue "E18" "  +1.18%  "
Set-TextValue "E19" "  -0.51%  "
Set-TextValue "D20" "12.73"
Set-TextValue "E20" "  -0.64%  "
Set-TextValue "E21" "  -0.49%  "
Set-TextValue "D22" "338.66"
Set-TextValue "E22" "  +2.03%  "
Set-TextValue "E23" "  +0.09%  "
Set-TextValue "D24" "0.511"
Set-TextValue "E24" "  +2.03%  "
Set-TextValue "D25" "66.46"
Set-TextValue "E25" "  +1.66%  "
Set-TextValue "D26" "0.169"
Set-TextValue "E26" "  -1.60%  "
Set-TextValue "E27" "  +0.11%  "
Set-TextValue "D28" "0.0₃0910"
Set-TextValue "E28" "  +0.74%  "
Set-TextValue "E29" "  -0.01%  "
Set-TextValue "E30" "  +0.06%  "
Set-TextValue "D31" "7.15"
Set-TextValue "E31" "  -1.54%  "
Set-TextValue "E32" "  +2.05%  "
Set-TextValue "D33" "20.85"
Set-TextValue "E33" "  +0.72%  "
Set-TextValue "E34" "  -1.47%  "
Set-TextValue "D35" "155.90"
Set-TextValue "E35" "  +0.83%  "
Set-TextValue "E36" "  +1.65%  "
Set-TextValue "E37" "  +1.49%  "
Set-TextValue "D38" "27.12"
Set-TextValue "E38" "  +0.36%  "
Set-TextValue "E39" "  -2.16%  "
Set-TextValue "D40" "0.0658"
Set-TextValue "E40" "  -3.02%  "
Set-TextValue "E41" "  +10.90%  "
Set-TextValue "E42" "  +0.11%  "
Set-TextValue "D43" "3.133.46"
Set-TextValue "E43" "  +0.98%  "
Set-TextValue "E44" "  +3.95%  "
Set-TextValue "E45" "  +0.48%  "
Set-TextValue "D47" "2.296.82"
Set-TextValue "E47" "  +0.91%  "
Set-TextValue "D48" "0.0258"
Set-TextValue "E48" "  -0.35%  "
Set-TextValue "D49" "0.974"
Set-TextValue "E49" "  +3.84%  "
Set-TextValue "E50" "  -1.35%  "
Set-TextValue "E51" "  +1.17%  "
